$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header styling updates (row 1): font name Arial -> Lato, fill color indexed 48 -> 42 ---
$header = $ws.Range("A1:B1")
$header.Font.Name = "Lato"
$header.Interior.ColorIndex = 35

# --- Extend the data rows, copying the existing row-3 formatting down to the new rows ---
$ws.Range("A3:B3").Copy()
$ws.Range("A4:B6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# --- Fill in the word/count data (climate keeps its row, new words added) ---
$ws.Range("A3").Value = "tundra"
$ws.Range("B3").Value = 1
$ws.Range("A4").Value = "temperature"
$ws.Range("B4").Value = 4
$ws.Range("A5").Value = "climate"
$ws.Range("B5").Value = 6
$ws.Range("A6").Value = "animals"
$ws.Range("B6").Value = 4
